# Applies the cryptos.xlsx "cryptos list" refresh (prices / 1h volume %
# updates, plus a RenderToken/Maker row swap in rows 45-46) as described
# by the commit "Updated cryptos list ... with GitHub Actions".
#
# Each target cell keeps its original text (inlineStr) semantics. For
# cells whose new text is also a valid numeric literal (e.g. "41.15"),
# we force a text number format before assigning the value and then
# reset the cell style back to "Normal" so no stray formatting/style is
# introduced, matching the original workbook's plain cell styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.191.76"
$ws.Range("D3").Value = "1.900.16"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.695"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +4.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.74%  "
$ws.Range("D14").Value = "2.176.44"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.738"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").Value = "1.895.66"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "35.180.52"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "242.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").Value = "4.128.54"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +18.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +17.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -13.06%  "
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0646"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.52%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.329.34"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.88%  "
